# Apply updated cryptocurrency price/volume figures (and the
# ImmutableX / RenzoRestakedETH row-order swap) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (cells in this sheet store numeric-
# looking prices as text, e.g. "66.351.25" / "0.999"). Forcing the
# number format to text before the write keeps Excel from silently
# re-typing these as floating point numbers, then the style is reset
# back to the original "Normal" so no visible formatting changes.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "66.351.25"
$ws.Range("E2").Value = "  -5.40%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.461.31"
$ws.Range("E3").Value = "  -6.58%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue $ws.Range("D5") "601.35"
$ws.Range("E5").Value = "  -7.64%  "

# Row 6
Set-TextValue $ws.Range("D6") "147.60"
$ws.Range("E6").Value = "  -9.32%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.461.00"
$ws.Range("E7").Value = "  -6.52%  "

# Row 8
$ws.Range("E8").Value = "  +0.19%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.477"
$ws.Range("E9").Value = "  -5.58%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.137"
$ws.Range("E10").Value = "  -7.42%  "

# Row 11
$ws.Range("E11").Value = "  -5.23%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.420"
$ws.Range("E12").Value = "  -6.38%  "

# Row 13
$ws.Range("E13").Value = "  -8.10%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.043.90"
$ws.Range("E14").Value = "  -6.33%  "

# Row 15
Set-TextValue $ws.Range("D15") "31.15"
$ws.Range("E15").Value = "  -5.79%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.445.34"
$ws.Range("E16").Value = "  -6.62%  "

# Row 17
Set-TextValue $ws.Range("D17") "66.378.70"
$ws.Range("E17").Value = "  -5.26%  "

# Row 18
$ws.Range("E18").Value = "  -0.47%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.31"
$ws.Range("E19").Value = "  -3.99%  "

# Row 20
Set-TextValue $ws.Range("D20") "14.91"
$ws.Range("E20").Value = "  -8.17%  "

# Row 21
Set-TextValue $ws.Range("D21") "439.48"
$ws.Range("E21").Value = "  -7.30%  "

# Row 22
Set-TextValue $ws.Range("D22") "8.92"
$ws.Range("E22").Value = "  -16.19%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.618"
$ws.Range("E23").Value = "  -5.74%  "

# Row 24
Set-TextValue $ws.Range("D24") "76.24"
$ws.Range("E24").Value = "  -5.00%  "

# Row 25
Set-TextValue $ws.Range("D25") "0.999"
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
Set-TextValue $ws.Range("D26") "3.603.21"
$ws.Range("E26").Value = "  -6.31%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.0000120"
$ws.Range("E27").Value = "  -7.95%  "

# Row 28
Set-TextValue $ws.Range("D28") "9.98"
$ws.Range("E28").Value = "  -10.20%  "

# Row 29
Set-TextValue $ws.Range("D29") "8.12"
$ws.Range("E29").Value = "  -12.68%  "

# Row 30
$ws.Range("E30").Value = "  -6.96%  "

# Row 31
$ws.Range("E31").Value = "  -10.66%  "

# Row 32
$ws.Range("E32").Value = "  +0.12%  "

# Row 33
Set-TextValue $ws.Range("D33") "25.34"
$ws.Range("E33").Value = "  -6.32%  "

# Row 34
$ws.Range("E34").Value = "  -6.82%  "

# Row 35
Set-TextValue $ws.Range("D35") "6.08"
$ws.Range("E35").Value = "  -8.22%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "1.83"
$ws.Range("E36").Value = "  -9.98%  "

# Row 37
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Range("D37") "3.449.30"
$ws.Range("E37").Value = "  -6.71%  "

# Row 38
Set-TextValue $ws.Range("D38") "7.87"
$ws.Range("E38").Value = "  -7.53%  "

# Row 39
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("E40").Value = "  +0.26%  "

# Row 41
Set-TextValue $ws.Range("D41") "172.05"
$ws.Range("E41").Value = "  -4.71%  "

# Row 42
$ws.Range("E42").Value = "  -7.21%  "

# Row 43
Set-TextValue $ws.Range("D43") "5.42"
$ws.Range("E43").Value = "  -8.90%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0852"
$ws.Range("E44").Value = "  -6.57%  "

# Row 45
$ws.Range("E45").Value = "  -6.13%  "

# Row 46
Set-TextValue $ws.Range("D46") "44.93"
$ws.Range("E46").Value = "  -4.72%  "

# Row 47
Set-TextValue $ws.Range("D47") "26.78"
$ws.Range("E47").Value = "  -8.86%  "

# Row 48
$ws.Range("E48").Value = "  -6.60%  "

# Row 49
Set-TextValue $ws.Range("D49") "7.47"
$ws.Range("E49").Value = "  -5.44%  "

# Row 50
Set-TextValue $ws.Range("D50") "2.43"
$ws.Range("E50").Value = "  -15.93%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.992"
$ws.Range("E51").Value = "  -7.61%  "
